$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save off the "before" values for the cells that rotate between rows 2, 3 and 4.
$A2 = $ws.Range("A2").Value2
$A3 = $ws.Range("A3").Value2
$A4 = $ws.Range("A4").Value2

$I2 = $ws.Range("I2").Value2
$I3 = $ws.Range("I3").Value2
$I4 = $ws.Range("I4").Value2

$Q2 = $ws.Range("Q2").Value2
$Q3 = $ws.Range("Q3").Value2
$Q4 = $ws.Range("Q4").Value2

$R2 = $ws.Range("R2").Value2
$R3 = $ws.Range("R3").Value2
$R4 = $ws.Range("R4").Value2

$AC3 = $ws.Range("AC3").Value2

# Apply the rotation: old row2 data -> row3, old row3 data -> row4, old row4 data -> row2
$ws.Range("A2").Value = $A4
$ws.Range("A3").Value = $A2
$ws.Range("A4").Value = $A3

$ws.Range("I2").Value = $I4
$ws.Range("I3").Value = $I2
$ws.Range("I4").Value = $I3

$ws.Range("Q2").Value = $Q4
$ws.Range("Q3").Value = $Q2
$ws.Range("Q4").Value = $Q3

$ws.Range("R2").Value = $R4
$ws.Range("R3").Value = $R2
$ws.Range("R4").Value = $R3

# The "i grävd grop" comment follows the row-3 data to row 4; row 3 no longer has it.
$ws.Range("AC3").ClearContents()
$ws.Range("AC4").Value = $AC3
